# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the
# 6135c9e7-7e16-4819-8d49-20d1572f7e07 file is now "Ready for handoff"
# (instead of "Handed back: in sync with en-US"), refreshes the related
# timestamps, records an Error Detail message on the target-language
# sheets, and widens the "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c1111b4de758a143abfd6751fb2be7904f5d2b4e/e2e/6135c9e7-7e16-4819-8d49-20d1572f7e07.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7057b79ee16901157278ef2851819ceaf3614bb1/e2e/6135c9e7-7e16-4819-8d49-20d1572f7e07.md."

# --- "Overview" sheet: row 3 is the 6135c9e7-... file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-06 21:06:35"

# --- "zh-cn" sheet: row 3 is the 6135c9e7-... file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-09-06 21:06:30"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- "de-de" sheet: row 3 is the 6135c9e7-... file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-09-06 21:06:35"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
